$wb = $excel.ActiveWorkbook

# --- Saturday Morning sheet ---
$ws1 = $wb.Worksheets.Item("Saturday Morning")
$ws1.Range("D30").Value = "Base Objective Function"
$ws1.Range("F30").Value = "Base Avg Transfer Wait Time"
$ws1.Range("D31").Value = 256520
$ws1.Range("F31").Formula = "=D31/D25"

# --- Saturday Evening sheet ---
$ws2 = $wb.Worksheets.Item("Saturday Evening")
$ws2.Range("D30").Value = "Base Objective Function"
$ws2.Range("F30").Value = "Base Avg Transfer Wait Time"
$ws2.Range("D31").Value = 228052
$ws2.Range("F31").Formula = "=D31/D25"

# Update selections / view
$ws2.Range("F31").Select()
$ws1.Range("G29").Select()
